$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Fix the wording in F6 (Agencia Hondureña de Aeronáutica Civil comunicado):
#    remove "por siete dias " from the sentence.
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = "En marco de la declaracion de emergencia nacional y las ultimas disposiciones emitidas en el Decreto Ejecutivo PCM -021-2020 en el que se decreta la suspensión  de algunas garantias constitucionales. Se notifica que se dispone el cierre de operaciones de los aeropuertos internacionales."

# ---------------------------------------------------------------------------
# 2) Fill in the previously-empty "Descarga Link" / "Descripción información"
#    cells for row 8 (Secretaría de Desarrollo Económico).
# ---------------------------------------------------------------------------
$ws.Range("F8").Value = "El gobierno de la República, a través de la Secretaria de Desarollo Económico, en el marco de la emergencia nacional sanitaria decretada, informa: Para garantizar que no se comentan abusus contra la población, a nivel nacional se decreta congelamiento de preciosos absoluto de los productos de la canasta básica y productos de higiene personal y de hogar."

# ---------------------------------------------------------------------------
# 3) Populate the new data row (row 9) - Instituto Nacional de Migración.
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "Instituto Nacional de Migración "
$ws.Range("C9").Value = "Trabajo"
$ws.Range("D9").Value = "Instituto Nacional de Migración tiene como misión fundamental ejercer el control y regulación como máxima autoridad en materia migratoria a nacionales y extranjeros en el marco de la protección de sus derechos y seguridad, en aplicación de la Ley de Migración y Extranjería y la Política Migratoria del Gobierno de la República, mediante una gestión migratoria moderna, dinámica y transparente.`n "
$ws.Range("F9").Value = "Siguiendo las disposiciones emitidas por el gobierno de Honduras el INM informa que los puntos de control migratorio se mantienen operando para: 1) Ingreso al país de hondureños residentes y diplomáticos. 2) Salida del terrotoio nacional a extranjeros. 3) Entrada y salida de transportistas de carga unicamente."
$ws.Range("I9").Value = "21/3/2020"
$ws.Range("J9").Value = "Honduras"

# ---------------------------------------------------------------------------
# 4) Hyperlinks for the new/changed URL cells (E8, E9, G9). Add() re-applies
#    hyperlink formatting with a brand new style, so re-stamp the intended
#    existing style (copied from a sibling cell) right after.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E8"), "https://covid19honduras.org/?q=Comunicados&page=8") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E9"), "https://covid19honduras.org/?q=Comunicados&page=7") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G9"), "https://covid19honduras.org/?q=flujos-migratorios") | Out-Null

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null

$ws.Range("G8").Copy() | Out-Null
$ws.Range("G9").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) Row 9 grows to the same height as row 8 now that it holds a full record.
# ---------------------------------------------------------------------------
$ws.Rows("9:9").RowHeight = 120

# ---------------------------------------------------------------------------
# 6) Update the saved view position / selection.
# ---------------------------------------------------------------------------
$ws.Range("D11").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 3
